# Update gh-pages output (generated at 456a3b4)
# Adds the new "南宁·花海演绎二次元水上派对" event (2024-10-05) to the
# "展览" and "全部类型" sheets, bumping every later row down by one and
# refreshing the "想去人数" (want-to-go) counters that changed between
# scrapes.
#
# NOTE: this runtime's `.Value` getter is unreliable (it echoes back the
# property signature instead of the cell's data), so every *read* below
# goes through `.Value2` instead; writes use `.Value2` too for symmetry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) -- new event lands at row 5
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Refresh counters on existing rows 2-4 (unrelated rows, values only)
$ws1.Range("F2").Value2 = 89
$ws1.Range("F3").Value2 = 371
$ws1.Range("F4").Value2 = 4798

# Row 6 (was row 5): push the old "熊喵M" row down by one.
$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value2 = $ws1.Range("B5").Value2
$ws1.Range("C6").Value2 = $ws1.Range("C5").Value2
$ws1.Range("D6").Value2 = $ws1.Range("D5").Value2
$ws1.Range("E6").Value2 = $ws1.Range("E5").Value2
$ws1.Range("F6").Value2 = $ws1.Range("F5").Value2
$ws1.Range("G6").Value2 = $ws1.Range("G5").Value2
$ws1.Range("H6").Value2 = $ws1.Range("H5").Value2
$ws1.Range("I6").Value2 = $ws1.Range("I5").Value2

# Row 7 (brand new row): the old "万圣漫控嘉年华10" row, with its refreshed count.
$ws1.Range("A6").Copy()
$ws1.Range("A7").PasteSpecial(-4122)
$ws1.Range("A7").Value2 = 6
$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value2 = "2024-11-02"
$ws1.Range("C7").Value2 = "南宁·万圣漫控嘉年华10"
$ws1.Range("D7").Value2 = "亭洪路45号 百益上河城"
$ws1.Range("E7").Value2 = "2024.11.02 11:00-11.03 22:00"
$ws1.Range("F7").Value2 = 483
$ws1.Range("G7").Value2 = 50
$ws1.Range("H7").Value2 = "https://show.bilibili.com/platform/detail.html?id=87820"
$ws1.Range("I7").Value2 = "//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg"

# Row 5: overwrite with the new event (do this last -- rows 6/7 above were
# copied from row 5's original contents).
$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value2 = "2024-10-05"
$ws1.Range("C5").Value2 = "南宁·花海演绎二次元水上派对"
$ws1.Range("D5").Value2 = "月湾路凤岭儿童公园北门东侧约70米 凤岭儿童公园"
$ws1.Range("E5").Value2 = "2024.10.05 14:00-10.05 21:00"
$ws1.Range("F5").Value2 = 10
$ws1.Range("G5").Value2 = 72
$ws1.Range("H5").Value2 = "https://show.bilibili.com/platform/detail.html?id=92559"
$ws1.Range("I5").Value2 = "//i0.hdslb.com/bfs/openplatform/202409/MTS1nDly1726642819177.png"

# ---------------------------------------------------------------------
# Sheet "全部类型" (all categories) -- new event lands at row 6
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

# Refresh counters on existing rows 2-4 (unrelated rows, values only)
$ws4.Range("F2").Value2 = 89
$ws4.Range("F3").Value2 = 371
$ws4.Range("F4").Value2 = 4798

# Row 9 (brand new row): the old "万圣漫控嘉年华10" row (was row 8), refreshed count.
$ws4.Range("A8").Copy()
$ws4.Range("A9").PasteSpecial(-4122)
$ws4.Range("A9").Value2 = 8
$ws4.Range("B9").NumberFormat = "@"
$ws4.Range("B9").Value2 = $ws4.Range("B8").Value2
$ws4.Range("C9").Value2 = $ws4.Range("C8").Value2
$ws4.Range("D9").Value2 = $ws4.Range("D8").Value2
$ws4.Range("E9").Value2 = $ws4.Range("E8").Value2
$ws4.Range("F9").Value2 = 483
$ws4.Range("G9").Value2 = $ws4.Range("G8").Value2
$ws4.Range("H9").Value2 = $ws4.Range("H8").Value2
$ws4.Range("I9").Value2 = $ws4.Range("I8").Value2

# Row 8 (was row 7): the old "熊喵M" row, pushed down by one.
$ws4.Range("B8").NumberFormat = "@"
$ws4.Range("B8").Value2 = $ws4.Range("B7").Value2
$ws4.Range("C8").Value2 = $ws4.Range("C7").Value2
$ws4.Range("D8").Value2 = $ws4.Range("D7").Value2
$ws4.Range("E8").Value2 = $ws4.Range("E7").Value2
$ws4.Range("F8").Value2 = $ws4.Range("F7").Value2
$ws4.Range("G8").Value2 = $ws4.Range("G7").Value2
$ws4.Range("H8").Value2 = $ws4.Range("H7").Value2
$ws4.Range("I8").Value2 = $ws4.Range("I7").Value2

# Row 7 (was row 6): the old "井草圣二" row, pushed down by one.
$ws4.Range("B7").NumberFormat = "@"
$ws4.Range("B7").Value2 = $ws4.Range("B6").Value2
$ws4.Range("C7").Value2 = $ws4.Range("C6").Value2
$ws4.Range("D7").Value2 = $ws4.Range("D6").Value2
$ws4.Range("E7").Value2 = $ws4.Range("E6").Value2
$ws4.Range("F7").Value2 = $ws4.Range("F6").Value2
$ws4.Range("G7").Value2 = $ws4.Range("G6").Value2
$ws4.Range("H7").Value2 = $ws4.Range("H6").Value2
$ws4.Range("I7").Value2 = $ws4.Range("I6").Value2

# Row 6: overwrite with the new event (do this last -- row 7 above was
# copied from row 6's original contents).
$ws4.Range("B6").NumberFormat = "@"
$ws4.Range("B6").Value2 = "2024-10-05"
$ws4.Range("C6").Value2 = "南宁·花海演绎二次元水上派对"
$ws4.Range("D6").Value2 = "月湾路凤岭儿童公园北门东侧约70米 凤岭儿童公园"
$ws4.Range("E6").Value2 = "2024.10.05 14:00-10.05 21:00"
$ws4.Range("F6").Value2 = 10
$ws4.Range("G6").Value2 = 72
$ws4.Range("H6").Value2 = "https://show.bilibili.com/platform/detail.html?id=92559"
$ws4.Range("I6").Value2 = "//i0.hdslb.com/bfs/openplatform/202409/MTS1nDly1726642819177.png"
